$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.060.09'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '3.773.94'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'628.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.38%  '
$ws.Range("D6").Value = "'165.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("D7").Value = '3.773.84'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D10").Value = "'0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("D12").Value = "'6.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").Value = "'35.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").Value = '4.408.31'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '3.753.90'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '69.097.49'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = "'17.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.19%  '
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").Value = "'7.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").Value = "'467.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = "'9.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").Value = "'82.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("E27").Value = '  +3.75%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '3.924.09'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").Value = "'7.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").Value = "'28.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").Value = "'0.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +18.85%  '
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '3.725.45'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("E40").Value = '  +3.38%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").Value = "'0.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = "'154.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("D46").Value = "'43.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = "'46.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("E49").Value = '  +4.10%  '
$ws.Range("E51").Value = '  -1.28%  '
